$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set C4 (Fim column) with the new date value, formatted as short date (mm-dd-yy, numFmtId 14)
$ws.Range("C4").Value = 42975.875
$ws.Range("C4").NumberFormat = "mm-dd-yy"

# Widen column C (Fim) to fit the newly added date; leave column B untouched
# so it keeps its original bestFit width (this naturally splits the merged
# "B:C" column-width entry into separate B and C entries).
$ws.Columns.Item(3).ColumnWidth = 26.6667

# Move the active selection from B4 to B5
$ws.Range("B5").Select()
